$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest scraped
# crypto data. Percent-change strings in column E are never auto-parsed as
# numbers (they contain a trailing "%" with padding spaces), but some of the
# "Price" strings in column D look like plain numbers (e.g. "210.74") - for
# those, force the cell to Text format first so Excel stores/keeps the exact
# original text (with its trailing zeros / dot-as-thousands-separator) instead
# of silently converting it to a numeric value.

$ws.Range("D2").Value = "28.316.43"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.560.82"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.74"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.34"
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.57"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "1.783.90"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "1.557.74"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "28.306.80"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.65"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.98"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.54"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +1.40%  "
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.20"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.87"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.05"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "1.380.07"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  +1.92%  "
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.518"
$ws.Range("E41").Value = "  -3.45%  "
$ws.Range("E42").Value = "  +2.59%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0471"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.780"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.09"
$ws.Range("E47").Value = "  -1.19%  "
$ws.Range("E48").Value = "  -6.31%  "
$ws.Range("D49").Value = "1.696.43"
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "85.33"
$ws.Range("E51").Value = "  -2.05%  "
